$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dcn"
$ws.Range("C2").Value = "Tlr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.15163
$ws.Range("H2").Value = 9.454890000000001
$ws.Range("I2").Value = 0.0006291248881010851
$ws.Range("J2").Value = 0.0006291248881010851
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 51.23401333333334
$ws.Range("N2").Value = 153.70204
$ws.Range("O2").Value = 0.9688226788583661
$ws.Range("P2").Value = 0.9688226788583661
$ws.Range("Q2").Value = 161.4706534417334
$ws.Range("R2").Value = 1453.2358809756
$ws.Range("S2").Value = 0.000609510459426563
$ws.Range("T2").Value = 0.000609510459426563

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dcn"
$ws.Range("C3").Value = "Tlr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.15163
$ws.Range("H3").Value = 9.454890000000001
$ws.Range("I3").Value = 0.0006291248881010851
$ws.Range("J3").Value = 0.0006291248881010851
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.648742666666666
$ws.Range("N3").Value = 4.946228
$ws.Range("O3").Value = 0.03117732114163389
$ws.Range("P3").Value = 0.03117732114163389
$ws.Range("Q3").Value = 5.196226850546666
$ws.Range("R3").Value = 46.76604165492
$ws.Range("S3").Value = [double]"1.961442867452202E-05"
$ws.Range("T3").Value = [double]"1.961442867452202E-05"

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Dcn"
$ws.Range("C4").Value = "Tlr2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4971.754394666666
$ws.Range("H4").Value = 14915.263184
$ws.Range("I4").Value = 0.9924561027819714
$ws.Range("J4").Value = 0.9924561027819713
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 51.23401333333334
$ws.Range("N4").Value = 153.70204
$ws.Range("O4").Value = 0.9688226788583661
$ws.Range("P4").Value = 0.9688226788583661
$ws.Range("Q4").Value = 254722.9309464106
$ws.Range("R4").Value = 2292506.378517695
$ws.Range("S4").Value = 0.9615139801465635
$ws.Range("T4").Value = 0.9615139801465634

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Dcn"
$ws.Range("C5").Value = "Tlr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4971.754394666666
$ws.Range("H5").Value = 14915.263184
$ws.Range("I5").Value = 0.9924561027819714
$ws.Range("J5").Value = 0.9924561027819713
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.648742666666666
$ws.Range("N5").Value = 4.946228
$ws.Range("O5").Value = 0.03117732114163389
$ws.Range("P5").Value = 0.03117732114163389
$ws.Range("Q5").Value = 8197.143598674438
$ws.Range("R5").Value = 73774.29238806994
$ws.Range("S5").Value = 0.03094212263540793
$ws.Range("T5").Value = 0.03094212263540793

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Dcn"
$ws.Range("C6").Value = "Tlr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 34.63986933333333
$ws.Range("H6").Value = 103.919608
$ws.Range("I6").Value = 0.006914772329927541
$ws.Range("J6").Value = 0.006914772329927542
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 51.23401333333334
$ws.Range("N6").Value = 153.70204
$ws.Range("O6").Value = 0.9688226788583661
$ws.Range("P6").Value = 0.9688226788583661
$ws.Range("Q6").Value = 1774.739527288924
$ws.Range("R6").Value = 15972.65574560032
$ws.Range("S6").Value = 0.006699188252376107
$ws.Range("T6").Value = 0.006699188252376108

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Dcn"
$ws.Range("C7").Value = "Tlr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 34.63986933333333
$ws.Range("H7").Value = 103.919608
$ws.Range("I7").Value = 0.006914772329927541
$ws.Range("J7").Value = 0.006914772329927542
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 1.648742666666666
$ws.Range("N7").Value = 4.946228
$ws.Range("O7").Value = 0.03117732114163389
$ws.Range("P7").Value = 0.03117732114163389
$ws.Range("Q7").Value = 57.11223053762487
$ws.Range("R7").Value = 514.010074838624
$ws.Range("S7").Value = 0.000215584077551435
$ws.Range("T7").Value = 0.000215584077551435

"Done"
